$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New order line items are stored as text (matching the existing rows),
# so force a text number format on the numeric-looking columns before
# writing the values to avoid Excel auto-converting them to numbers.
$ws.Range("C3:E4").NumberFormat = "@"

$ws.Range("A3").Value = "104APRICOT28"
$ws.Range("B3").Value = "Apricot - Dried"
$ws.Range("C3").Value = "2"
$ws.Range("D3").Value = "128.49"
$ws.Range("E3").Value = "256.98"

$ws.Range("A4").Value = "760RD5"
$ws.Range("B4").Value = "Container - Alur Deli (5oz)"
$ws.Range("C4").Value = "1"
$ws.Range("D4").Value = "97.99"
$ws.Range("E4").Value = "97.99"
